$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L2, N2 (plain values)
$ws.Range("L2").Value = 157
$ws.Range("N2").Value = 100

# Update R2, T2 (plain values)
$ws.Range("R2").Value = 332
$ws.Range("T2").Value = 64

# New row 3 totals (M3, S3)
$ws.Range("M3").Value = 531
$ws.Range("S3").Value = 724

# M2, S2 become formulas referencing the new row 3 totals
$ws.Range("M2").Formula = "=M3-(L2+N2)"
$ws.Range("S2").Formula = "=S3-(R2+T2)"

# M3/S3 keep the "data" style (fill+border+center, like L2/N2/R2/T2 ...)
$ws.Range("L2").Copy()
$ws.Range("M3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# M2/S2 switch to the plain "center" style (same as B2) now that they hold formulas
$ws.Range("B2").Copy()
$ws.Range("M2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Clear out the old row 5 (A5 = 160) which is no longer part of the data
$ws.Range("A5").ClearContents()

# Update the active selection to match the target view state
$ws.Range("L10").Select()
